{"js": "// Update each two-digit multiplication-problem cell in the practice-sheet\n// table (\"<operand>\u00d7<operand>=\") to the new operand pair from the commit's\n// diff. All 100 problems are unique strings, so a targeted search+replace\n// per pair (in document order) reproduces the edit exactly, touching only\n// the <w:t> text of each run and leaving every other run/paragraph/table\n// property untouched.\nconst pairs = [\n  [\"13\u00d728=\", \"26\u00d759=\"],\n  [\"42\u00d778=\", \"36\u00d710=\"],\n  [\"18\u00d790=\", \"93\u00d713=\"],\n  [\"54\u00d743=\", \"16\u00d766=\"],\n  [\"16\u00d744=\", \"65\u00d724=\"],\n  [\"83\u00d745=\", \"72\u00d776=\"],\n  [\"31\u00d728=\", \"44\u00d737=\"],\n  [\"35\u00d715=\", \"58\u00d774=\"],\n  [\"23\u00d718=\", \"14\u00d768=\"],\n  [\"66\u00d724=\", \"60\u00d776=\"],\n  [\"71\u00d787=\", \"31\u00d714=\"],\n  [\"25\u00d726=\", \"19\u00d717=\"],\n  [\"46\u00d757=\", \"63\u00d783=\"],\n  [\"40\u00d712=\", \"35\u00d711=\"],\n  [\"93\u00d720=\", \"41\u00d762=\"],\n  [\"99\u00d710=\", \"47\u00d728=\"],\n  [\"45\u00d732=\", \"11\u00d754=\"],\n  [\"79\u00d763=\", \"86\u00d780=\"],\n  [\"100\u00d749=\", \"30\u00d780=\"],\n  [\"90\u00d757=\", \"82\u00d713=\"],\n  [\"58\u00d712=\", \"36\u00d794=\"],\n  [\"31\u00d748=\", \"85\u00d792=\"],\n  [\"13\u00d735=\", \"78\u00d761=\"],\n  [\"31\u00d751=\", \"98\u00d754=\"],\n  [\"25\u00d786=\", \"13\u00d739=\"],\n  [\"12\u00d731=\", \"29\u00d776=\"],\n  [\"46\u00d715=\", \"27\u00d774=\"],\n  [\"16\u00d792=\", \"17\u00d729=\"],\n  [\"59\u00d759=\", \"36\u00d777=\"],\n  [\"57\u00d793=\", \"79\u00d760=\"],\n  [\"76\u00d759=\", \"23\u00d779=\"],\n  [\"83\u00d717=\", \"26\u00d781=\"],\n  [\"16\u00d746=\", \"20\u00d749=\"],\n  [\"57\u00d794=\", \"17\u00d718=\"],\n  [\"84\u00d742=\", \"13\u00d749=\"],\n  [\"63\u00d715=\", \"33\u00d795=\"],\n  [\"59\u00d746=\", \"51\u00d774=\"],\n  [\"44\u00d747=\", \"19\u00d794=\"],\n  [\"22\u00d742=\", \"58\u00d729=\"],\n  [\"29\u00d764=\", \"94\u00d717=\"],\n  [\"57\u00d736=\", \"26\u00d791=\"],\n  [\"69\u00d727=\", \"10\u00d797=\"],\n  [\"77\u00d751=\", \"66\u00d717=\"],\n  [\"43\u00d712=\", \"56\u00d794=\"],\n  [\"70\u00d736=\", \"53\u00d721=\"],\n  [\"50\u00d798=\", \"61\u00d769=\"],\n  [\"96\u00d721=\", \"44\u00d713=\"],\n  [\"96\u00d793=\", \"31\u00d743=\"],\n  [\"52\u00d738=\", \"82\u00d727=\"],\n  [\"72\u00d746=\", \"56\u00d766=\"],\n  [\"74\u00d773=\", \"47\u00d712=\"],\n  [\"88\u00d738=\", \"90\u00d745=\"],\n  [\"49\u00d711=\", \"40\u00d721=\"],\n  [\"16\u00d724=\", \"80\u00d731=\"],\n  [\"99\u00d738=\", \"68\u00d725=\"],\n  [\"17\u00d799=\", \"91\u00d757=\"],\n  [\"68\u00d757=\", \"96\u00d796=\"],\n  [\"16\u00d715=\", \"53\u00d748=\"],\n  [\"14\u00d714=\", \"84\u00d739=\"],\n  [\"54\u00d789=\", \"82\u00d781=\"],\n  [\"50\u00d748=\", \"83\u00d778=\"],\n  [\"66\u00d718=\", \"74\u00d766=\"],\n  [\"21\u00d799=\", \"45\u00d756=\"],\n  [\"85\u00d772=\", \"73\u00d740=\"],\n  [\"17\u00d715=\", \"99\u00d725=\"],\n  [\"51\u00d739=\", \"89\u00d736=\"],\n  [\"28\u00d756=\", \"24\u00d798=\"],\n  [\"74\u00d781=\", \"85\u00d797=\"],\n  [\"73\u00d755=\", \"35\u00d785=\"],\n  [\"66\u00d731=\", \"34\u00d750=\"],\n  [\"68\u00d745=\", \"15\u00d719=\"],\n  [\"65\u00d747=\", \"61\u00d774=\"],\n  [\"67\u00d747=\", \"78\u00d737=\"],\n  [\"66\u00d799=\", \"60\u00d742=\"],\n  [\"82\u00d772=\", \"75\u00d786=\"],\n  [\"94\u00d749=\", \"80\u00d766=\"],\n  [\"86\u00d727=\", \"48\u00d786=\"],\n  [\"80\u00d797=\", \"30\u00d728=\"],\n  [\"13\u00d756=\", \"25\u00d778=\"],\n  [\"47\u00d763=\", \"39\u00d731=\"],\n  [\"76\u00d777=\", \"36\u00d772=\"],\n  [\"10\u00d782=\", \"60\u00d783=\"],\n  [\"80\u00d764=\", \"17\u00d746=\"],\n  [\"95\u00d789=\", \"50\u00d789=\"],\n  [\"55\u00d787=\", \"34\u00d737=\"],\n  [\"78\u00d718=\", \"62\u00d775=\"],\n  [\"67\u00d723=\", \"81\u00d711=\"],\n  [\"50\u00d769=\", \"55\u00d741=\"],\n  [\"98\u00d731=\", \"49\u00d739=\"],\n  [\"87\u00d710=\", \"53\u00d7100=\"],\n  [\"20\u00d720=\", \"10\u00d751=\"],\n  [\"63\u00d737=\", \"24\u00d717=\"],\n  [\"90\u00d795=\", \"33\u00d720=\"],\n  [\"65\u00d757=\", \"87\u00d763=\"],\n  [\"14\u00d790=\", \"62\u00d726=\"],\n  [\"15\u00d736=\", \"35\u00d762=\"],\n  [\"75\u00d740=\", \"23\u00d772=\"],\n  [\"99\u00d723=\", \"55\u00d789=\"],\n  [\"28\u00d768=\", \"16\u00d785=\"],\n  [\"32\u00d735=\", \"97\u00d754=\"]\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly 1 match for \"${oldText}\", found ${results.items.length}`);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update each two-digit multiplication-problem cell in the practice-sheet\n# table (\"<operand>x<operand>=\") to the new operand pair from the commit's\n# diff. All 100 problems are unique strings, so a targeted Find/Replace per\n# pair (in document order) reproduces the edit exactly, touching only the\n# <w:t> text of each run and leaving every other run/paragraph/table\n# property untouched.\n$d = $word.ActiveDocument\n$pairs = @(\n  @('13\u00d728=', '26\u00d759='),\n  @('42\u00d778=', '36\u00d710='),\n  @('18\u00d790=', '93\u00d713='),\n  @('54\u00d743=', '16\u00d766='),\n  @('16\u00d744=', '65\u00d724='),\n  @('83\u00d745=', '72\u00d776='),\n  @('31\u00d728=', '44\u00d737='),\n  @('35\u00d715=', '58\u00d774='),\n  @('23\u00d718=', '14\u00d768='),\n  @('66\u00d724=', '60\u00d776='),\n  @('71\u00d787=', '31\u00d714='),\n  @('25\u00d726=', '19\u00d717='),\n  @('46\u00d757=', '63\u00d783='),\n  @('40\u00d712=', '35\u00d711='),\n  @('93\u00d720=', '41\u00d762='),\n  @('99\u00d710=', '47\u00d728='),\n  @('45\u00d732=', '11\u00d754='),\n  @('79\u00d763=', '86\u00d780='),\n  @('100\u00d749=', '30\u00d780='),\n  @('90\u00d757=', '82\u00d713='),\n  @('58\u00d712=', '36\u00d794='),\n  @('31\u00d748=', '85\u00d792='),\n  @('13\u00d735=', '78\u00d761='),\n  @('31\u00d751=', '98\u00d754='),\n  @('25\u00d786=', '13\u00d739='),\n  @('12\u00d731=', '29\u00d776='),\n  @('46\u00d715=', '27\u00d774='),\n  @('16\u00d792=', '17\u00d729='),\n  @('59\u00d759=', '36\u00d777='),\n  @('57\u00d793=', '79\u00d760='),\n  @('76\u00d759=', '23\u00d779='),\n  @('83\u00d717=', '26\u00d781='),\n  @('16\u00d746=', '20\u00d749='),\n  @('57\u00d794=', '17\u00d718='),\n  @('84\u00d742=', '13\u00d749='),\n  @('63\u00d715=', '33\u00d795='),\n  @('59\u00d746=', '51\u00d774='),\n  @('44\u00d747=', '19\u00d794='),\n  @('22\u00d742=', '58\u00d729='),\n  @('29\u00d764=', '94\u00d717='),\n  @('57\u00d736=', '26\u00d791='),\n  @('69\u00d727=', '10\u00d797='),\n  @('77\u00d751=', '66\u00d717='),\n  @('43\u00d712=', '56\u00d794='),\n  @('70\u00d736=', '53\u00d721='),\n  @('50\u00d798=', '61\u00d769='),\n  @('96\u00d721=', '44\u00d713='),\n  @('96\u00d793=', '31\u00d743='),\n  @('52\u00d738=', '82\u00d727='),\n  @('72\u00d746=', '56\u00d766='),\n  @('74\u00d773=', '47\u00d712='),\n  @('88\u00d738=', '90\u00d745='),\n  @('49\u00d711=', '40\u00d721='),\n  @('16\u00d724=', '80\u00d731='),\n  @('99\u00d738=', '68\u00d725='),\n  @('17\u00d799=', '91\u00d757='),\n  @('68\u00d757=', '96\u00d796='),\n  @('16\u00d715=', '53\u00d748='),\n  @('14\u00d714=', '84\u00d739='),\n  @('54\u00d789=', '82\u00d781='),\n  @('50\u00d748=', '83\u00d778='),\n  @('66\u00d718=', '74\u00d766='),\n  @('21\u00d799=', '45\u00d756='),\n  @('85\u00d772=', '73\u00d740='),\n  @('17\u00d715=', '99\u00d725='),\n  @('51\u00d739=', '89\u00d736='),\n  @('28\u00d756=', '24\u00d798='),\n  @('74\u00d781=', '85\u00d797='),\n  @('73\u00d755=', '35\u00d785='),\n  @('66\u00d731=', '34\u00d750='),\n  @('68\u00d745=', '15\u00d719='),\n  @('65\u00d747=', '61\u00d774='),\n  @('67\u00d747=', '78\u00d737='),\n  @('66\u00d799=', '60\u00d742='),\n  @('82\u00d772=', '75\u00d786='),\n  @('94\u00d749=', '80\u00d766='),\n  @('86\u00d727=', '48\u00d786='),\n  @('80\u00d797=', '30\u00d728='),\n  @('13\u00d756=', '25\u00d778='),\n  @('47\u00d763=', '39\u00d731='),\n  @('76\u00d777=', '36\u00d772='),\n  @('10\u00d782=', '60\u00d783='),\n  @('80\u00d764=', '17\u00d746='),\n  @('95\u00d789=', '50\u00d789='),\n  @('55\u00d787=', '34\u00d737='),\n  @('78\u00d718=', '62\u00d775='),\n  @('67\u00d723=', '81\u00d711='),\n  @('50\u00d769=', '55\u00d741='),\n  @('98\u00d731=', '49\u00d739='),\n  @('87\u00d710=', '53\u00d7100='),\n  @('20\u00d720=', '10\u00d751='),\n  @('63\u00d737=', '24\u00d717='),\n  @('90\u00d795=', '33\u00d720='),\n  @('65\u00d757=', '87\u00d763='),\n  @('14\u00d790=', '62\u00d726='),\n  @('15\u00d736=', '35\u00d762='),\n  @('75\u00d740=', '23\u00d772='),\n  @('99\u00d723=', '55\u00d789='),\n  @('28\u00d768=', '16\u00d785='),\n  @('32\u00d735=', '97\u00d754=')\n)\n\nforeach ($p in $pairs) {\n  $old = $p[0]\n  $new = $p[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $ok = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n  if (-not $ok) {\n    throw \"Replacement failed for '$old' -> '$new'\"\n  }\n}\n"}
